# Update committee and affiliations data on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G13 gained a value of 1 (loc_extended/spc column for an existing committee member)
$ws.Range("G13").Value = 1

# New committee member added as row 48: first, last, institution, and spc = 1
$ws.Range("A48").Value = "Marci"
$ws.Range("B48").Value = "Rückbeil"
$ws.Range("C48").Value = "Sanofi"
$ws.Range("G48").Value = 1

# Match the author's final selection/view position
$ws.Range("G49").Select()
